$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$col = $ws.Columns.Item(7)  # Column G ("Recorded By")

$replaced = 0
$firstFound = $col.Find($target, [Type]::Missing, [Type]::Missing, 1)
if ($firstFound) {
    $firstAddr = $firstFound.Address()
    $cell = $firstFound
    do {
        $cell.Value2 = $replacement
        $replaced++
        $cell = $col.FindNext($cell)
    } while ($cell -and $cell.Address() -ne $firstAddr)
}

Write-Host "Replaced $replaced occurrences of '$target' with '$replacement' in column G."
